# Apply the two kinds of changes described by the diff:
#  1. Refresh the cached text of the "datetimeFigureOut" auto-date field
#     (on every slide layout, the slide master, and the notes master)
#     from 10/6/2023 -> 10/24/2023.
#  2. Re-crop / reposition the picture on slide 27 (bottom-crop it and
#     move/resize the frame accordingly).

$p = $ppt.ActivePresentation

$newDate = "10/24/2023"

# --- 1. Date placeholder on each slide layout -----------------------------
$m = $p.SlideMaster

# Map: layout index (1-based, matches slideLayoutN.xml) -> shape index of
# the shape that holds the <a:fld type="datetimeFigureOut"> date field.
$layoutDateShape = @{
    1  = 3
    2  = 3
    3  = 3
    4  = 4
    5  = 6
    6  = 2
    7  = 1
    8  = 4
    9  = 4
    10 = 3
    11 = 3
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    $shIdx = $layoutDateShape[$li]
    if ($shIdx) {
        $shape = $layout.Shapes.Item($shIdx)
        $shape.TextFrame.TextRange.Text = $newDate
    }
}

# --- Date placeholder on the slide master ----------------------------------
$m.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# --- Date placeholder on the notes master -----------------------------------
$nm = $p.NotesMaster
$nm.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# --- 2. Picture crop / resize on slide 27 -----------------------------------
$slide = $p.Slides.Item(27)
$pic = $slide.Shapes.Item(2)

# Crop off the bottom ~35.385% of the source image (matches <a:srcRect b="35385"/>).
$pic.PictureFormat.CropBottom = 49.36247977233413

# New frame position/size (EMU 1339849,2675255 / 8704457x2059306), expressed
# in points (1 pt = 12700 EMU) and nudged to account for the host's f32
# rounding so the serialized EMU values land exactly on target.
$pic.Left = 105.49992370605469
$pic.Top = 210.65000915527344
$pic.Width = 685.3903198242188
$pic.Height = 162.15008544921875
